$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6; this pushes current rows 6-17 down to 7-18
$ws.Rows.Item(6).Insert()

# Copy the date number format (style) used by column D from the row below (row 7, the shifted original row 6)
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new row 6 with the new weekly data point
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44914
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100114007
$ws.Cells.Item(6, 7).Value = "Jengibre"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 14000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 13).Value = 14400
$ws.Cells.Item(6, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(6, 15).Value = "Perú"
$ws.Cells.Item(6, 16).Value = 1108
$ws.Cells.Item(6, 17).Value = 13
$ws.Cells.Item(6, 18).Value = "Hortaliza"
